$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Date placeholder fields: "8/20/2020" -> "8/21/2020"
#    These live on the slide master, every slide layout, and the notes
#    master (one "Date Placeholder" shape each).
# ---------------------------------------------------------------------------
function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "8/20/2020") {
                $tr.Text = "8/21/2020"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShapes $master.Shapes

$layouts = $master.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DateShapes $layouts.Item($j).Shapes
}

$notesMaster = $p.NotesMaster
Update-DateShapes $notesMaster.Shapes

# ---------------------------------------------------------------------------
# 2) Slide 26 title: "Structural Induction on the Binary Heap"
#    -> "Mathematical Induction on the Height"
# ---------------------------------------------------------------------------
$s26 = $p.Slides.Item(26)
$s26.Shapes.Item(1).TextFrame.TextRange.Text = "Mathematical Induction on the Height"

# ---------------------------------------------------------------------------
# 3) Slide 39 (Summary) body text: collapse the three runs describing
#    "structural induction" into a single sentence about discrete
#    structures, leaving the rest of the paragraph (line breaks + the
#    "We will cover Part II ..." runs) untouched.
# ---------------------------------------------------------------------------
$s39 = $p.Slides.Item(39)
$body = $s39.Shapes.Item(2)
$tr39 = $body.TextFrame.TextRange

$fullText = $tr39.Text
$oldPhrase = "We have also illustrated a new variant of mathematical induction known as structural induction. "
$startIdx = $fullText.IndexOf($oldPhrase)
if ($startIdx -ge 0) {
    $sub = $tr39.Characters($startIdx + 1, $oldPhrase.Length)
    $sub.Text = "We have also illustrated an application of mathematical induction on discrete structures."
}
